$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell without letting Excel's
# "looks like a number/date" auto-detection convert it to a numeric/date
# cell (and without leaving any NumberFormat residue on the cell's style).
# We do this by computing the literal string via TEXT(), then collapsing
# the formula down to its static result with a values-only paste.
function Set-TextValue {
    param($range, [string]$text)
    $escaped = $text.Replace('"', '""')
    $range.Formula = '=TEXT("' + $escaped + '","@")'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

# New "Vigencia" header in F1, matching the style used by the other header cells
$ws.Range("A1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "Vigencia"

# Row 2 (new record inserted at top): Roberto  Laforcada / IRSA / Maximiliano Godoy
Set-TextValue $ws.Range("A2") "33013638"
Set-TextValue $ws.Range("B2") "Roberto  Laforcada"
Set-TextValue $ws.Range("C2") "IRSA"
Set-TextValue $ws.Range("D2") "Maximiliano Godoy"
Set-TextValue $ws.Range("E2") "2025-11-14"

# Row 3: Analia Belen Miño, now under IRSA
Set-TextValue $ws.Range("A3") "35115887"
Set-TextValue $ws.Range("B3") "Analia Belen Miño"
Set-TextValue $ws.Range("C3") "IRSA"
Set-TextValue $ws.Range("D3") "Roberto Laforcada"
Set-TextValue $ws.Range("E3") "2025-11-14"

# Row 4: Paris Laforcada, now under IRSA
Set-TextValue $ws.Range("A4") "53412356"
Set-TextValue $ws.Range("B4") "Paris Laforcada"
Set-TextValue $ws.Range("C4") "IRSA"
Set-TextValue $ws.Range("D4") "Roberto Laforcada"
Set-TextValue $ws.Range("E4") "2025-11-14"

# Row 5: Patrick Laforcada, now under IRSA, plus a Vigencia date value
Set-TextValue $ws.Range("A5") "59610581"
Set-TextValue $ws.Range("B5") "Patrick Laforcada"
Set-TextValue $ws.Range("C5") "IRSA"
Set-TextValue $ws.Range("D5") "Roberto Laforcada"
Set-TextValue $ws.Range("E5") "2025-11-14"

$ws.Range("F5").Value = 46003
$ws.Range("F5").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("F5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
